$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5370729
$ws.Range("J17").Value = 5668997
$ws.Range("L17").Value = 17006991
$ws.Range("N17").Value = -17007327

$ws.Range("H43").Value = 595.9167
$ws.Range("I43").Value = 310
$ws.Range("K43").Value = 310
$ws.Range("M43").Value = -241

$ws.Range("H64").Value = 4100
$ws.Range("J64").Value = 4700
$ws.Range("L64").Value = 4700
$ws.Range("N64").Value = -5196

$ws.Range("H67").Value = 4100
$ws.Range("J67").Value = 4700
$ws.Range("L67").Value = 4700
$ws.Range("N67").Value = -6416

$ws.Range("H113").Value = 100005550
$ws.Range("I113").Value = 200003700
$ws.Range("K113").Value = 200003700
$ws.Range("M113").Value = -200000446

$ws.Range("H129").Value = 295267
$ws.Range("J129").Value = 295267
$ws.Range("L129").Value = 885801
$ws.Range("N129").Value = -895801

$ws.Range("H132").Value = 3867.9614
$ws.Range("I132").Value = 4464.8096
$ws.Range("K132").Value = 13394.4288
$ws.Range("M132").Value = -10864.4288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4540.14
$ws.Range("I32").Value = 3712.314
$ws.Range("J32").Value = 9625.357
$ws.Range("K32").Value = 3712.314
$ws.Range("L32").Value = 9625.357
$ws.Range("M32").Value = -3425.314
$ws.Range("N32").Value = -10199.357

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H110").Value = 973.5454999999999
$ws.Range("I110").Value = 878.8889
$ws.Range("K110").Value = 878.8889
$ws.Range("M110").Value = 1166.1111

$ws.Range("H122").Value = 2138.8125
$ws.Range("I122").Value = 1665.9286
$ws.Range("J122").Value = 5449
$ws.Range("K122").Value = 4997.7858
$ws.Range("L122").Value = 16347
$ws.Range("M122").Value = -2547.7858
$ws.Range("N122").Value = -21247

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2083
$ws.Range("I20").Value = 2367.9285
$ws.Range("J20").Value = 1639.7778
$ws.Range("K20").Value = 2367.9285
$ws.Range("L20").Value = 1639.7778
$ws.Range("M20").Value = -2120.9285
$ws.Range("N20").Value = -2133.7778

$ws.Range("H86").Value = 1708.1842
$ws.Range("I86").Value = 1581.64
$ws.Range("J86").Value = 1951.5385
$ws.Range("K86").Value = 1581.64
$ws.Range("L86").Value = 1951.5385
$ws.Range("M86").Value = -458.6400000000001
$ws.Range("N86").Value = -4197.538500000001

$ws.Range("H89").Value = 1708.1842
$ws.Range("I89").Value = 1581.64
$ws.Range("J89").Value = 1951.5385
$ws.Range("K89").Value = 7908.200000000001
$ws.Range("L89").Value = 9757.692500000001
$ws.Range("M89").Value = -2292.200000000001
$ws.Range("N89").Value = -20989.6925

$ws.Range("H134").Value = 4217.1665
$ws.Range("I134").Value = 4579
$ws.Range("J134").Value = 3028.2856
$ws.Range("K134").Value = 13737
$ws.Range("L134").Value = 9084.856800000001
$ws.Range("M134").Value = -11202
$ws.Range("N134").Value = -14154.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4941.7744
$ws.Range("I31").Value = 3989.818
$ws.Range("J31").Value = 5465.35
$ws.Range("K31").Value = 3989.818
$ws.Range("L31").Value = 5465.35
$ws.Range("M31").Value = -3694.818
$ws.Range("N31").Value = -6055.35

$ws.Range("H34").Value = 4941.7744
$ws.Range("I34").Value = 3989.818
$ws.Range("J34").Value = 5465.35
$ws.Range("K34").Value = 3989.818
$ws.Range("L34").Value = 5465.35
$ws.Range("M34").Value = -3787.818
$ws.Range("N34").Value = -5869.35

$ws.Range("H105").Value = 1053.36
$ws.Range("J105").Value = 1589.6
$ws.Range("L105").Value = 1589.6
$ws.Range("N105").Value = -5083.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 52.61111
$ws.Range("I12").Value = 8.625
$ws.Range("J12").Value = 87.8
$ws.Range("K12").Value = 25.875
$ws.Range("L12").Value = 263.4
$ws.Range("M12").Value = 147.125
$ws.Range("N12").Value = -609.4

$ws.Range("H69").Value = 1971.1428
$ws.Range("J69").Value = 1971.1428
$ws.Range("L69").Value = 5913.428400000001
$ws.Range("N69").Value = -7535.428400000001

$ws.Range("H72").Value = 1971.1428
$ws.Range("J72").Value = 1971.1428
$ws.Range("L72").Value = 17740.2852
$ws.Range("N72").Value = -25852.2852

$ws.Range("H131").Value = 698.7217000000001
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 698.7217000000001
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2096.1651
$ws.Range("N131").Value = -12176.1651
$ws.Range("M131").ClearContents()

$ws.Range("H138").Value = 2509.5264
$ws.Range("I138").Value = 2015.3846
$ws.Range("K138").Value = 6046.1538
$ws.Range("M138").Value = -906.1538

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1413.8462
$ws.Range("I97").Value = 1475.5555
$ws.Range("J97").Value = 1275
$ws.Range("K97").Value = 1475.5555
$ws.Range("L97").Value = 1275
$ws.Range("M97").Value = -979.5554999999999
$ws.Range("N97").Value = -2267

$ws.Range("H132").Value = 23408.154
$ws.Range("I132").Value = 4480
$ws.Range("J132").Value = 86502
$ws.Range("K132").Value = 13440
$ws.Range("L132").Value = 259506
$ws.Range("M132").Value = -10910
$ws.Range("N132").Value = -264566

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3219.8667
$ws.Range("I40").Value = 2914.1428
$ws.Range("J40").Value = 7500
$ws.Range("K40").Value = 2914.1428
$ws.Range("L40").Value = 7500
$ws.Range("M40").Value = -2778.1428
$ws.Range("N40").Value = -7772

$ws.Range("H46").Value = 1910.2778
$ws.Range("I46").Value = 2268.077
$ws.Range("K46").Value = 2268.077
$ws.Range("M46").Value = -2080.077

$ws.Range("H104").Value = 20607.715
$ws.Range("J104").Value = 20607.715
$ws.Range("L104").Value = 20607.715
$ws.Range("N104").Value = -27595.715

$ws.Range("H122").Value = 983677
$ws.Range("J122").Value = 4498.2856
$ws.Range("L122").Value = 13494.8568
$ws.Range("N122").Value = -18394.8568

$ws.Range("H132").Value = 1006881.2
$ws.Range("I132").Value = 2009682.6
$ws.Range("J132").Value = 4079.6667
$ws.Range("K132").Value = 6029047.800000001
$ws.Range("L132").Value = 12239.0001
$ws.Range("M132").Value = -6026517.800000001
$ws.Range("N132").Value = -17299.0001

$ws.Range("H136").Value = 1293.0968
$ws.Range("I136").Value = 1188.3704
$ws.Range("K136").Value = 3565.1112
$ws.Range("M136").Value = -1015.1112

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 12475.375
$ws.Range("J4").Value = 12475.375
$ws.Range("L4").Value = 12475.375
$ws.Range("N4").Value = -12701.375

$ws.Range("H62").Value = 4605.7
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 4757.125
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 4757.125
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -6005.125

$ws.Range("H65").Value = 4605.7
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 4757.125
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 23785.625
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -30025.625
